$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.757.51"
$ws.Range("E2").Value = "  +1.10%  "

$ws.Range("D3").Value = "3.159.00"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "530.42"
$ws.Range("E5").Value = "  -0.81%  "

$ws.Range("D6").Value = "139.74"
$ws.Range("E6").Value = "  +0.67%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "0.539"
$ws.Range("E8").Value = "  +15.26%  "

$ws.Range("E9").Value = "  +0.54%  "

$ws.Range("D10").Value = "0.437"
$ws.Range("E10").Value = "  +5.52%  "

$ws.Range("D11").Value = "0.112"
$ws.Range("E11").Value = "  +4.10%  "

$ws.Range("E12").Value = "  +3.01%  "

$ws.Range("D13").Value = "3.700.69"
$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("E14").Value = "  +1.09%  "

$ws.Range("E15").Value = "  +4.97%  "

$ws.Range("D16").Value = "58.776.99"
$ws.Range("E16").Value = "  +0.99%  "

$ws.Range("D17").Value = "6.26"
$ws.Range("E17").Value = "  +3.59%  "

$ws.Range("D18").Value = "3.142.22"
$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("D19").Value = "'13.00"
$ws.Range("E19").Value = "  +2.56%  "

$ws.Range("D20").Value = "8.15"
$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").Value = "372.13"
$ws.Range("E21").Value = "  +3.44%  "

$ws.Range("E22").Value = "  +1.43%  "

$ws.Range("E23").Value = "  +0.31%  "

$ws.Range("D24").Value = "'0.530"
$ws.Range("E24").Value = "  +4.67%  "

$ws.Range("D25").Value = "69.71"
$ws.Range("E25").Value = "  +1.00%  "

$ws.Range("D26").Value = "0.168"
$ws.Range("E26").Value = "  +0.92%  "

$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("D28").Value = "8.29"
$ws.Range("E28").Value = "  +13.42%  "

$ws.Range("E29").Value = "  -1.11%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "22.28"
$ws.Range("E30").Value = "  +3.33%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.88"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").Value = "6.12"
$ws.Range("E32").Value = "  +0.35%  "

$ws.Range("E33").Value = "  +2.76%  "

$ws.Range("E34").Value = "  +1.05%  "

$ws.Range("D35").Value = "6.33"
$ws.Range("E35").Value = "  +4.21%  "

$ws.Range("D36").Value = "158.38"
$ws.Range("E36").Value = "  -0.80%  "

$ws.Range("E37").Value = "  +4.81%  "

$ws.Range("D38").Value = "25.14"
$ws.Range("E38").Value = "  -2.82%  "

$ws.Range("E39").Value = "  -0.36%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.0686"
$ws.Range("E40").Value = "  +2.36%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.634.90"
$ws.Range("E41").Value = "  +5.21%  "

$ws.Range("E42").Value = "  +6.30%  "

$ws.Range("D43").Value = "39.09"
$ws.Range("E43").Value = "  +4.56%  "

$ws.Range("D44").Value = "0.0286"
$ws.Range("E44").Value = "  +6.47%  "

$ws.Range("D45").Value = "0.709"
$ws.Range("E45").Value = "  +0.70%  "

$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").Value = "3.197.42"
$ws.Range("E47").Value = "  +0.43%  "

$ws.Range("D48").Value = "0.104"
$ws.Range("E48").Value = "  +13.26%  "

$ws.Range("E49").Value = "  +2.41%  "

$ws.Range("D50").Value = "'0.980"
$ws.Range("E50").Value = "  -0.59%  "

$ws.Range("E51").Value = "  +2.36%  "
